$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3-5 (the extra AWB rows), keep header + one data row
$ws.Rows.Item(3).Resize(3).Delete() | Out-Null

# Update remaining data row values
$ws.Range("A2").Value = 4714411638593
$ws.Range("B2").Value = 20

# Column width change for column A (target stored width 18.7109375)
$ws.Columns.Item(1).ColumnWidth = 17.8

# Update selection
$ws.Range("I13").Select() | Out-Null

# Page setup - portrait orientation
$ws.PageSetup.Orientation = 1
